$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1481  # F2: 1479 -> 1481
$ws.Cells.Item(3, 6).Value = 114  # F3: 113 -> 114
$ws.Cells.Item(5, 6).Value = 7392  # F5: 7386 -> 7392
$ws.Cells.Item(6, 6).Value = 74  # F6: 73 -> 74
$ws.Cells.Item(7, 6).Value = 4759  # F7: 4755 -> 4759
$ws.Cells.Item(8, 6).Value = 7004  # F8: 7001 -> 7004
$ws.Cells.Item(10, 6).Value = 271  # F10: 268 -> 271
$ws.Cells.Item(12, 6).Value = 855  # F12: 853 -> 855
$ws.Cells.Item(13, 6).Value = 166  # F13: 165 -> 166
$ws.Cells.Item(14, 6).Value = 50  # F14: 49 -> 50
$ws.Cells.Item(15, 6).Value = 1161  # F15: 1160 -> 1161
$ws.Cells.Item(17, 6).Value = 157  # F17: 156 -> 157
$ws.Cells.Item(19, 6).Value = 223  # F19: 221 -> 223
$ws.Cells.Item(21, 6).Value = 1154  # F21: 1152 -> 1154
$ws.Cells.Item(25, 6).Value = 1221  # F25: 1219 -> 1221
$ws.Cells.Item(26, 6).Value = 43  # F26: 42 -> 43
$ws.Cells.Item(30, 6).Value = 169  # F30: 166 -> 169
$ws.Cells.Item(33, 6).Value = 90  # F33: 87 -> 90
$ws.Cells.Item(36, 6).Value = 426  # F36: 425 -> 426
$ws.Cells.Item(37, 6).Value = 67  # F37: 66 -> 67
$ws.Cells.Item(39, 6).Value = 367  # F39: 366 -> 367
$ws.Cells.Item(40, 6).Value = 1197  # F40: 1196 -> 1197
$ws.Cells.Item(41, 6).Value = 575  # F41: 574 -> 575
$ws.Cells.Item(43, 6).Value = 20  # F43: 19 -> 20
$ws.Cells.Item(44, 6).Value = 18  # F44: 17 -> 18

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(26, 6).Value = 632  # F26: 631 -> 632
$ws.Cells.Item(28, 6).Value = 25  # F28: 24 -> 25
$ws.Cells.Item(31, 6).Value = 851  # F31: 847 -> 851
$ws.Cells.Item(33, 6).Value = 605  # F33: 604 -> 605
$ws.Cells.Item(39, 6).Value = 106  # F39: 105 -> 106
$ws.Cells.Item(43, 6).Value = 74  # F43: 73 -> 74

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 727  # F4: 726 -> 727
$ws.Cells.Item(6, 6).Value = 668  # F6: 666 -> 668
$ws.Cells.Item(8, 6).Value = 1564  # F8: 1555 -> 1564
$ws.Cells.Item(9, 6).Value = 2457  # F9: 2446 -> 2457

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 727  # F2: 726 -> 727
$ws.Cells.Item(3, 6).Value = 1481  # F3: 1479 -> 1481
$ws.Cells.Item(5, 6).Value = 114  # F5: 113 -> 114
$ws.Cells.Item(6, 6).Value = 668  # F6: 666 -> 668
$ws.Cells.Item(7, 6).Value = 668  # F7: 666 -> 668
$ws.Cells.Item(8, 6).Value = 7392  # F8: 7386 -> 7392
$ws.Cells.Item(9, 6).Value = 74  # F9: 73 -> 74
$ws.Cells.Item(10, 6).Value = 4759  # F10: 4755 -> 4759
$ws.Cells.Item(12, 6).Value = 7004  # F12: 7001 -> 7004
$ws.Cells.Item(13, 6).Value = 271  # F13: 268 -> 271
$ws.Cells.Item(16, 6).Value = 855  # F16: 853 -> 855
$ws.Cells.Item(17, 6).Value = 166  # F17: 165 -> 166
$ws.Cells.Item(18, 6).Value = 1564  # F18: 1555 -> 1564
$ws.Cells.Item(19, 6).Value = 2457  # F19: 2446 -> 2457
$ws.Cells.Item(21, 6).Value = 50  # F21: 49 -> 50
$ws.Cells.Item(22, 6).Value = 1161  # F22: 1160 -> 1161
$ws.Cells.Item(24, 6).Value = 157  # F24: 156 -> 157
$ws.Cells.Item(25, 6).Value = 223  # F25: 221 -> 223
$ws.Cells.Item(26, 6).Value = 1154  # F26: 1152 -> 1154
$ws.Cells.Item(27, 6).Value = 632  # F27: 631 -> 632
$ws.Cells.Item(29, 6).Value = 1221  # F29: 1219 -> 1221
$ws.Cells.Item(31, 6).Value = 169  # F31: 166 -> 169
$ws.Cells.Item(32, 6).Value = 25  # F32: 24 -> 25
$ws.Cells.Item(33, 6).Value = 851  # F33: 847 -> 851
$ws.Cells.Item(35, 6).Value = 90  # F35: 87 -> 90
$ws.Cells.Item(38, 6).Value = 605  # F38: 604 -> 605
$ws.Cells.Item(39, 6).Value = 67  # F39: 66 -> 67
$ws.Cells.Item(42, 6).Value = 367  # F42: 366 -> 367
$ws.Cells.Item(43, 6).Value = 575  # F43: 574 -> 575
$ws.Cells.Item(45, 6).Value = 106  # F45: 105 -> 106
$ws.Cells.Item(48, 6).Value = 74  # F48: 73 -> 74
$ws.Cells.Item(49, 6).Value = 20  # F49: 19 -> 20
$ws.Cells.Item(50, 6).Value = 18  # F50: 17 -> 18
